$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.406
$ws.Range("B3").Value = 6.701000000000001
$ws.Range("C5").Value = -12.799
$ws.Range("D5").Value = -8.309000000000001
$ws.Range("E7").Value = 13.078
$ws.Range("D9").Value = -7.634
$ws.Range("D11").Value = -8.284000000000001
$ws.Range("E11").Value = 12.87
$ws.Range("B14").Value = 6.331
$ws.Range("E19").Value = 12.756
$ws.Range("B21").Value = 6.434
$ws.Range("D21").Value = -7.7
$ws.Range("E21").Value = 12.054
$ws.Range("B23").Value = 6.547
$ws.Range("B25").Value = 6.103999999999999
